$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(2, 1).Value = "30 Oct 2025, 01:10 AM"

# --- Top Gainers sheet: new entry SKMEGGPROD pushes into row 36, ---
# --- shifting the existing rows 36-75 down to 37-76, and the former ---
# --- row 76 (CGPOWER) falls off the bottom of the list. ---
$ws = $wb.Worksheets.Item("Top Gainers")

# Insert a fresh row at 36, shifting rows 36..76 down to 37..77.
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the new top-of-list entry.
$ws.Cells.Item(36, 1).Value = "🚀"
$ws.Cells.Item(36, 2).Value = "SKMEGGPROD"
$ws.Cells.Item(36, 3).Value = 4.9959
$ws.Cells.Item(36, 4).Value = 6.6906
$ws.Cells.Item(36, 5).Value = 23.7638

# Drop the row that fell off the end of the table (old row 76, now 77).
$ws.Rows.Item(77).Delete()
